# Update the "Estado de Cuenta" workbook:
#  - Remove the three extra worker rows (Carlos Ruben Pajaro Medina, Ronald Jose
#    Morales Ricardo, Francisco Andres Martinez Aguirre), keeping only Kerlis
#    Paola Contreras Caña.
#  - Refresh the aggregate totals (Valor Mora, Cant. Trabajadores, Cant. Periodos)
#    to reflect the remaining single worker / single period.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Delete rows 17-19 (the three workers being removed), shifting rows 20+ (the
# signature block) up so it lands on rows 21/22 as in the target layout.
$ws.Range("B17:J19").EntireRow.Delete() | Out-Null

# Update the summary figures to match the now-single remaining worker/period.
$ws.Range("E11").Value = 12133
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Column D ("Nombre Trabajador") no longer needs to fit the longest removed
# name, so its best-fit width shrinks along with the data (closest width the
# COM ColumnWidth rounding allows to the target 31.1796875 stored units).
$ws.Columns.Item(4).ColumnWidth = 30.3
